$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order of data rows after sorting each year block so that
# Oct, Nov, Dec come first, followed by Jan..Sep (matches target diff).
$data = @(
    @("2014-10", 94.967, 96.9919),
    @("2014-11", 93.98390000000001, 96.8069),
    @("2014-12", 93.0339, 97.7856),
    @("2014-01", 94.2403, 92.626),
    @("2014-02", 93.54940000000001, 92.53019999999999),
    @("2014-03", 92.46510000000001, 91.313),
    @("2014-04", 94.0853, 91.1688),
    @("2014-05", 94.30240000000001, 92.10980000000001),
    @("2014-06", 95.6863, 93.8792),
    @("2014-07", 97.7135, 95.56959999999999),
    @("2014-08", 96.8207, 96.1876),
    @("2014-09", 96.0857, 96.7236),
    @("2015-10", 86.7, 103.3),
    @("2015-11", 86.91549999999999, 103.4774),
    @("2015-12", 86.7807, 103.1929),
    @("2015-01", 90.2513, 98.2157),
    @("2015-02", 89.98569999999999, 97.4384),
    @("2015-03", 91.1819, 97.0228),
    @("2015-04", 91.7997, 98.3883),
    @("2015-05", 93.45950000000001, 99.25409999999999),
    @("2015-06", 90.6116, 99.5181),
    @("2015-07", 87.9521, 99.5557),
    @("2015-08", 86.245, 100.7725),
    @("2015-09", 86.38509999999999, 101.9467),
    @("2016-10", 96.90000000000001, 105.8),
    @("2016-11", 99.2, 105),
    @("2016-12", 104.3, 105.9),
    @("2016-01", 88.2967, 100.094),
    @("2016-02", 88.7299, 100.5276),
    @("2016-03", 90.6583, 102.8247),
    @("2016-04", 90.90519999999999, 103.8232),
    @("2016-05", 89.59999999999999, 102.8),
    @("2016-06", 90.40000000000001, 102.7),
    @("2016-07", 91.8, 102.7),
    @("2016-08", 94.40000000000001, 104.4),
    @("2016-09", 95, 105.9),
    @("2017-10", 114, 101.1),
    @("2017-11", 113.5, 101),
    @("2017-12", 109.5, 100.3),
    @("2017-01", 109.3, 109.2),
    @("2017-02", 111.9, 111.7),
    @("2017-03", 110.4, 112.9),
    @("2017-04", 108.3, 110.7),
    @("2017-05", 107.2, 108.3),
    @("2017-06", 108.7, 105.8),
    @("2017-07", 110.4, 105),
    @("2017-08", 110.9, 103.1),
    @("2017-09", 113, 101.9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
